# Applies the tracked changes:
#   1. Removes the (default, no-op) shading from the six header-table
#      cells that carried an explicit <w:shd val="clear" color="auto"
#      fill="auto"/> by resetting Shading to its "no fill / automatic"
#      state through the Word object model.
#   2. Adds a new paragraph right after "Anexam prezentei cereri
#      documentatia aferenta." asking that the notice be picked up in
#      person from the PMI registry.

$d = $word.ActiveDocument

# --- 1. Clear the shading on the first table's cells -----------------
$tbl = $d.Tables.Item(1)
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $row = $tbl.Rows.Item($r)
    for ($c = 1; $c -le $row.Cells.Count; $c++) {
        $cell = $row.Cells.Item($c)
        $cell.Shading.Texture = 0
        $cell.Shading.ForegroundPatternColor = -16777216
        $cell.Shading.BackgroundPatternColor = -16777216
    }
}

# --- 2. Insert the new paragraph --------------------------------------
# Use a fresh Range (instead of $d.Paragraphs directly) to look up the
# paragraph list - keeps the paragraph numbering correctly anchored to
# the whole document even after the table cells above were touched.
$full = $d.Content
$target = $d.Content
$found = $target.Find.Execute(
    "Anexăm prezentei cereri documentaţia aferentă.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $anchorPara = -1
    for ($i = 1; $i -le $full.Paragraphs.Count; $i++) {
        $p = $full.Paragraphs.Item($i)
        if ($p.Range.Start -le $target.Start -and $p.Range.End -ge $target.End) {
            $anchorPara = $i
        }
    }

    if ($anchorPara -gt 0) {
        $src = $full.Paragraphs.Item($anchorPara)
        $src.Range.InsertParagraphAfter()
        $newFull = $d.Content
        $newPara = $newFull.Paragraphs.Item($anchorPara + 1)
        $newPara.Range.Text = "Menționez că doresc să ridic avizul de la registratura PMI."
    }
}
